# Update excess mortality plots
# Applies the data refresh to excess_mortality_provinces.xlsx:
#   - A handful of weekly province death counts (rows 125-150) were
#     corrected/updated; all ROUND(...) percentage-change formulas in
#     columns AE:AP recalculate automatically from these inputs.
#   - Row 151 ("2022 week 45") gets its comparison columns (N:AA, AC:AD)
#     populated plus the matching AE:AP percentage-change formulas, mirroring
#     the pattern already used by the rows above it.
#   - The active selection is moved to reflect where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 125 ---
$ws.Range("U125").Value = 415

# --- Row 130 ---
$ws.Range("S130").Value = 204

# --- Row 132 ---
$ws.Range("W132").Value = 481

# --- Row 139 ---
$ws.Range("P139").Value = 116

# --- Row 140 ---
$ws.Range("AA140").Value = 231

# --- Row 141 ---
$ws.Range("W141").Value = 413

# --- Row 142 ---
$ws.Range("P142").Value = 123
$ws.Range("X142").Value = 595
$ws.Range("Z142").Value = 452

# --- Row 143 ---
$ws.Range("X143").Value = 607

# --- Row 144 ---
$ws.Range("W144").Value = 489
$ws.Range("X144").Value = 607

# --- Row 145 ---
$ws.Range("R145").Value = 120
$ws.Range("Z145").Value = 432

# --- Row 146 ---
$ws.Range("T146").Value = 51
$ws.Range("W146").Value = 501
$ws.Range("X146").Value = 663

# --- Row 147 ---
$ws.Range("T147").Value = 50
$ws.Range("V147").Value = 217
$ws.Range("W147").Value = 476
$ws.Range("X147").Value = 690
$ws.Range("Y147").Value = 88
$ws.Range("Z147").Value = 486

# --- Row 148 ---
$ws.Range("S148").Value = 252
$ws.Range("W148").Value = 499
$ws.Range("X148").Value = 659
$ws.Range("AA148").Value = 228

# --- Row 149 ---
$ws.Range("U149").Value = 444
$ws.Range("W149").Value = 513
$ws.Range("X149").Value = 640
$ws.Range("AA149").Value = 249

# --- Row 150 ---
$ws.Range("Q150").Value = 134
$ws.Range("R150").Value = 133
$ws.Range("S150").Value = 227
$ws.Range("U150").Value = 388
$ws.Range("V150").Value = 233
$ws.Range("W150").Value = 472
$ws.Range("X150").Value = 652
$ws.Range("Y150").Value = 82
$ws.Range("Z150").Value = 453
$ws.Range("AA150").Value = 241

# --- Row 151 ("2022 week 45") : new comparison data + computed columns ---
$ws.Range("N151").Value = 2022
$ws.Range("O151").Value = 45
$ws.Range("P151").Value = 117
$ws.Range("Q151").Value = 132
$ws.Range("R151").Value = 135
$ws.Range("S151").Value = 222
$ws.Range("T151").Value = 63
$ws.Range("U151").Value = 431
$ws.Range("V151").Value = 235
$ws.Range("W151").Value = 470
$ws.Range("X151").Value = 646
$ws.Range("Y151").Value = 76
$ws.Range("Z151").Value = 493
$ws.Range("AA151").Value = 261
$ws.Range("AC151").Value = 2022
$ws.Range("AD151").Value = 45

$ws.Range("AE151").Formula = "=ROUND((P151-B151)/B151*100,2)"
$ws.Range("AF151").Formula = "=ROUND((Q151-C151)/C151*100,2)"
$ws.Range("AG151").Formula = "=ROUND((R151-D151)/D151*100,2)"
$ws.Range("AH151").Formula = "=ROUND((S151-E151)/E151*100,2)"
$ws.Range("AI151").Formula = "=ROUND((T151-F151)/F151*100,2)"
$ws.Range("AJ151").Formula = "=ROUND((U151-G151)/G151*100,2)"
$ws.Range("AK151").Formula = "=ROUND((V151-H151)/H151*100,2)"
$ws.Range("AL151").Formula = "=ROUND((W151-I151)/I151*100,2)"
$ws.Range("AM151").Formula = "=ROUND((X151-J151)/J151*100,2)"
$ws.Range("AN151").Formula = "=ROUND((Y151-K151)/K151*100,2)"
$ws.Range("AO151").Formula = "=ROUND((Z151-L151)/L151*100,2)"
$ws.Range("AP151").Formula = "=ROUND((AA151-M151)/M151*100,2)"

# --- Reflect the author's final on-screen selection/scroll position ---
$ws.Range("AJ120").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
